# Applies the "Updated symbol list" crypto-price refresh for Sat Dec 31
# 2022 20:00 UTC GitHub Actions run.
#
# For each touched row, updates Coin (B), Link (C), Price (D),
# Volume(1h) (E) and Hora (G) cells. The worksheet stores every value as
# literal text (no real numbers), so each numeric-looking cell has its
# NumberFormat forced to "@" (Text) immediately before the assignment;
# otherwise Excel would auto-coerce things like "20", "0.89%" or
# "0.06981" into real numbers/percentages and mangle the exact digits
# (leading/trailing zeros, "%" suffix, etc.) that the source feed uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; D = "246.92"; E = "0.89%"; G = "20" },
    @{ Row = 3; D = "26.43"; E = "5.29%"; G = "20" },
    @{ Row = 4; E = "1.94%"; G = "20" },
    @{ Row = 5; D = "0.05601"; E = "-0.26%"; G = "20" },
    @{ Row = 6; D = "6.493"; E = "-0.78%"; G = "20" },
    @{ Row = 7; D = "0.8131"; E = "0.48%"; G = "20" },
    @{ Row = 8; E = "0.71%"; G = "20" },
    @{ Row = 9; B = "MandalaExchangeToken"; C = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D = "0.06981"; E = "0.18%"; G = "20" },
    @{ Row = 10; B = "BitrueCoin"; C = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D = "0.02848"; E = "0.22%"; G = "20" },
    @{ Row = 11; B = "BitMartToken"; C = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D = "0.09403"; E = "-0.04%"; G = "20" },
    @{ Row = 12; B = "BitForexToken"; C = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D = "0.001511"; E = "-1.15%"; G = "20" },
    @{ Row = 13; B = "TigerCash"; C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D = "0.006227"; E = "1.99%"; G = "20" },
    @{ Row = 14; B = "LEO"; C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D = "3.611"; E = "3.18%"; G = "20" },
    @{ Row = 15; B = "GateToken"; C = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; D = "3.014"; E = "0.31%"; G = "20" },
    @{ Row = 16; B = "BTSEToken"; C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; D = "2.055"; E = "-1.73%"; G = "20" },
    @{ Row = 17; B = "One"; C = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; D = "0.0005969"; E = "-0.16%"; G = "20" },
    @{ Row = 18; B = "BitpandaEcosystemToken"; C = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; D = "0.3176"; E = "-0.69%"; G = "20" },
    @{ Row = 19; B = "WazirX"; C = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; D = "0.1336"; E = "-0.02%"; G = "20" },
    @{ Row = 20; D = "0.03175"; E = "-2.46%"; G = "20" },
    @{ Row = 21; D = "0.1297"; E = "0.48%"; G = "20" },
    @{ Row = 22; D = "3.737"; E = "-0.25%"; G = "20" },
    @{ Row = 23; E = "-0.74%"; G = "20" },
    @{ Row = 24; E = "-1.45%"; G = "20" },
    @{ Row = 25; D = "0.001247"; E = "0.41%"; G = "20" },
    @{ Row = 26; D = "0.004589"; E = "1.42%"; G = "20" },
    @{ Row = 27; D = "0.00009598"; E = "-1.12%"; G = "20" },
    @{ Row = 28; G = "20" },
    @{ Row = 29; G = "20" },
    @{ Row = 30; G = "20" },
    @{ Row = 31; G = "20" },
    @{ Row = 32; G = "20" },
    @{ Row = 33; G = "20" },
    @{ Row = 34; G = "20" },
    @{ Row = 35; G = "20" },
    @{ Row = 36; G = "20" },
    @{ Row = 37; G = "20" },
    @{ Row = 38; G = "20" },
    @{ Row = 39; G = "20" },
    @{ Row = 40; D = "0.03671"; E = "0.61%"; G = "20" },
    @{ Row = 41; D = "0.006179"; E = "83.63%"; G = "20" },
    @{ Row = 42; D = "0.1058"; E = "-22.58%"; G = "20" },
    @{ Row = 43; D = "0.002500"; E = "-8.25%"; G = "20" },
    @{ Row = 44; D = "0.008937"; E = "10.53%"; G = "20" },
    @{ Row = 45; D = "0.00005357"; E = "1.58%"; G = "20" },
    @{ Row = 46; E = "0.01%"; G = "20" },
    @{ Row = 47; G = "20" },
    @{ Row = 48; D = "0.002615"; E = "28.06%"; G = "20" },
    @{ Row = 49; D = "0.00002100"; E = "0.01%"; G = "20" },
    @{ Row = 50; D = "0.0002000"; E = "0.01%"; G = "20" },
    @{ Row = 51; G = "20" },
)

foreach ($item in $data) {
    $r = $item.Row

    if ($item.ContainsKey("B")) {
        $ws.Range("B$r").Value = $item.B
    }
    if ($item.ContainsKey("C")) {
        $ws.Range("C$r").Value = $item.C
    }
    if ($item.ContainsKey("D")) {
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Range("E$r").NumberFormat = "@"
        $ws.Range("E$r").Value = $item.E
    }
    if ($item.ContainsKey("G")) {
        $ws.Range("G$r").NumberFormat = "@"
        $ws.Range("G$r").Value = $item.G
    }
}
